$wb = $excel.ActiveWorkbook

# Rename "Sheet2" to "NinzaAutomation"
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Name = "NinzaAutomation"

# Clear old data in A1:B1
$ws.Range("A1:B1").ClearContents()

# Move the data to F6 (Price) and G6 (Product Name)
$ws.Range("F6").Value = "Price"
$ws.Range("G6").Value = "Product Name"
